$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (currently sitting alone in its own
#    empty paragraph, right after the Github link paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. In the "TEN/WWD/889" paragraph, delete the extra "W" so the text becomes
#    "TEN/WD/889", then re-insert the "_GoBack" bookmark (collapsed) right
#    between "TEN/W" and "D/889".
$d.Content.Find.Execute("TEN/WWD/889", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TEN/WD/889", 2)

$p8 = $d.Paragraphs(8)
$bmRange = $p8.Range.Duplicate
$bmRange.Collapse(1)
$bmRange.MoveStart(1, 5)
$d.Bookmarks.Add("_GoBack", $bmRange)
